$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.921.37"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3
$ws.Range("D3").Value = "'1.762.74"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4
$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  -0.64%  "

# Row 5
$ws.Range("D5").Value = "'321.37"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").Value = "'0.9982"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("D7").Value = "'0.4233"
$ws.Range("E7").Value = "  -3.60%  "

# Row 8
$ws.Range("D8").Value = "'0.3584"
$ws.Range("E8").Value = "  -2.92%  "

# Row 9
$ws.Range("D9").Value = "'44.15"
$ws.Range("E9").Value = "  -1.50%  "

# Row 10
$ws.Range("D10").Value = "'0.07434"
$ws.Range("E10").Value = "  -2.57%  "

# Row 11
$ws.Range("D11").Value = "'1.099"
$ws.Range("E11").Value = "  -1.07%  "

# Row 12
$ws.Range("D12").Value = "'0.9974"
$ws.Range("E12").Value = "  -0.59%  "

# Row 13
$ws.Range("D13").Value = "'21.50"
$ws.Range("E13").Value = "  -0.18%  "

# Row 14
$ws.Range("D14").Value = "'6.078"
$ws.Range("E14").Value = "  -0.98%  "

# Row 15
$ws.Range("D15").Value = "'7.293"
$ws.Range("E15").Value = "  -1.45%  "

# Row 16
$ws.Range("D16").Value = "'1.782.76"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17
$ws.Range("D17").Value = "'90.78"
$ws.Range("E17").Value = "  +0.79%  "

# Row 18
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("E18").Value = "  -1.11%  "

# Row 19
$ws.Range("D19").Value = "'0.06363"
$ws.Range("E19").Value = "  +1.95%  "

# Row 20
$ws.Range("D20").Value = "'0.9982"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("D21").Value = "'17.12"

# Row 22
$ws.Range("D22").Value = "'5.960"
$ws.Range("E22").Value = "  -3.35%  "

# Row 23
$ws.Range("D23").Value = "'27.903.16"
$ws.Range("E23").Value = "  +1.40%  "

# Row 24
$ws.Range("D24").Value = "'11.25"
$ws.Range("E24").Value = "  -2.20%  "

# Row 25
$ws.Range("D25").Value = "'2.142"
$ws.Range("E25").Value = "  -7.03%  "

# Row 26
$ws.Range("D26").Value = "'159.73"
$ws.Range("E26").Value = "  +4.49%  "

# Row 27
$ws.Range("D27").Value = "'20.09"
$ws.Range("E27").Value = "  -1.83%  "

# Row 28
$ws.Range("D28").Value = "'1.981.72"
$ws.Range("E28").Value = "  +1.23%  "

# Row 29
$ws.Range("D29").Value = "'2.133"
$ws.Range("E29").Value = "  -7.19%  "

# Row 30
$ws.Range("D30").Value = "'125.08"
$ws.Range("E30").Value = "  -1.83%  "

# Row 31
$ws.Range("D31").Value = "'1.164"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32
$ws.Range("D32").Value = "'5.664"

# Row 33
$ws.Range("D33").Value = "'0.09036"
$ws.Range("E33").Value = "  -1.43%  "

# Row 34
$ws.Range("D34").Value = "'3.507"
$ws.Range("E34").Value = "  -3.47%  "

# Row 35
$ws.Range("D35").Value = "'12.64"
$ws.Range("E35").Value = "  +0.50%  "

# Row 36
$ws.Range("D36").Value = "'0.02305"
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.038"
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06058"
$ws.Range("E38").Value = "  -0.70%  "

# Row 39
$ws.Range("D39").Value = "'0.2098"
$ws.Range("E39").Value = "  -2.31%  "

# Row 40
$ws.Range("D40").Value = "'0.6371"
$ws.Range("E40").Value = "  -0.67%  "

# Row 41
$ws.Range("E41").Value = "  +0.83%  "

# Row 42
$ws.Range("D42").Value = "'0.9977"
$ws.Range("E42").Value = "  -0.35%  "

# Row 43
$ws.Range("D43").Value = "'7.874"
$ws.Range("E43").Value = "  -0.79%  "

# Row 44
$ws.Range("D44").Value = "'1.397"
$ws.Range("E44").Value = "  +0.42%  "

# Row 45
$ws.Range("D45").Value = "'13.54"
$ws.Range("E45").Value = "  -0.97%  "

# Row 46
$ws.Range("D46").Value = "'0.5942"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").Value = "'3.691"
$ws.Range("E47").Value = "  -0.64%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.987"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'123.31"
$ws.Range("E49").Value = "  -2.13%  "

# Row 50
$ws.Range("D50").Value = "'1.149"
$ws.Range("E50").Value = "  +1.51%  "

# Row 51
$ws.Range("D51").Value = "'0.06861"
$ws.Range("E51").Value = "  -0.17%  "
